# Update final entrants in each challenge
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 26250
    3  = 15248
    4  = 15842
    5  = 14872
    6  = 10861
    7  = 10795
    8  = 10742
    9  = 10860
    10 = 8176
    11 = 8732
    12 = 8042
    13 = 6348
    14 = 6161
    15 = 5632
    16 = 5371
    17 = 5192
    18 = 4614
    19 = 4373
    20 = 4088
    21 = 4125
    22 = 4002
    23 = 3806
    24 = 3694
    25 = 3300
    26 = 3230
    27 = 3850
    28 = 3187
    29 = 2768
    30 = 2307
    31 = 2109
    33 = 1734
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 2).Value = $values[$row]
}

# D33 previously held the shared string "@9:17am"; clear it so the now-unused
# string gets dropped from sharedStrings.xml.
$ws.Range("D33").ClearContents()

# Update the active selection to match the recorded cursor position.
$ws.Range("B32").Select()
